$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.693.32"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.45"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.58"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6523"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07492"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2971"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.50"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07640"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.51"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.045"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6863"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.47"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009578"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.128"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.728.50"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.106.25"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.12"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.61"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.718"
$ws.Range("E23").Value = "  +4.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.21"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1422"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.529"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.84"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06043"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.271"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.140"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.066"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.869"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.183"
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7263"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.599"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.801"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01787"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.201.30"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.284"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9121"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.019.15"
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.19"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.55"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.342"
$ws.Range("E47").Value = "  +9.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000121"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4049"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.159"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("E51").Value = "  +3.30%  "
